# Update the Rules workbook to reflect the "GIT UPDATE" commit:
#  - Cell E8 ("Good Morning") becomes "GIT UPDATE"
#  - Selection moves to E8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"

$ws.Range("E8").Select()
